$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date-formatted cell so the new rows reuse the existing "date" cell
# style instead of Excel synthesizing a brand-new number format.
$ws.Range("C35").Copy($ws.Range("C36"))
$ws.Range("C35").Copy($ws.Range("C37"))

# --- Row 36 : case #35 ---
$ws.Range("A36").Value = 35
$ws.Range("B36").Value = "msg"
$ws.Range("C36").Value = 44383
$ws.Range("D36").Value = "MCAST"
$ws.Range("E36").Value = "lookalike"
$ws.Range("F36").Value = "change"
$ws.Range("G36").Value = "en"
$ws.Range("H36").Value = "no"
$ws.Range("I36").Value = "Confirm number to continue using service"
$ws.Range("J36").Value = "BOV"

# --- Row 37 : case #36 ---
$ws.Range("A37").Value = 36
$ws.Range("B37").Value = "msg"
$ws.Range("C37").Value = 44383
$ws.Range("D37").Value = "MCAST"
$ws.Range("E37").Value = "lookalike"
$ws.Range("F37").Value = "delivery"
$ws.Range("G37").Value = "mt"
$ws.Range("H37").Value = "no"
$ws.Range("I37").Value = "confirm address for postal delivery"
$ws.Range("J37").Value = "DHL"

# Move the view the way the author left it: active selection on K36.
$ws.Activate()
$ws.Range("K36").Select()
